$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 75 values with revised figures
$ws.Range("B75").Value = 65605
$ws.Range("C75").Value = 9547
$ws.Range("D75").Value = 7419
$ws.Range("I75").Value = 7541
$ws.Range("K75").Value = 7068
$ws.Range("L75").Value = 17086
$ws.Range("M75").Value = 14605
$ws.Range("R75").Value = 19200
$ws.Range("S75").Value = -21045
$ws.Range("T75").Value = 86649
$ws.Range("Z75").Value = 1422
$ws.Range("AB75").Value = 1422
$ws.Range("AC75").Value = 12153

# Add new row 76 for the 01-04-2021 quarter
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").Style = "Normal"
$ws.Range("B76").Value = 66979
$ws.Range("C76").Value = 10464
$ws.Range("D76").Value = 8078
$ws.Range("E76").Value = 2387
$ws.Range("F76").Value = 9582
$ws.Range("G76").Value = 564
$ws.Range("H76").Value = 9018
$ws.Range("I76").Value = 8916
$ws.Range("J76").Value = 1157
$ws.Range("K76").Value = 7759
$ws.Range("L76").Value = 16791
$ws.Range("M76").Value = 14478
$ws.Range("N76").Value = 2313
$ws.Range("O76").Value = 0
$ws.Range("P76").Value = 30
$ws.Range("Q76").Value = 30
$ws.Range("R76").Value = 21195
$ws.Range("S76").Value = -20558
$ws.Range("T76").Value = 87537
$ws.Range("U76").Value = 0
$ws.Range("V76").Value = 0
$ws.Range("W76").Value = 72568
$ws.Range("X76").Value = 3703
$ws.Range("Y76").Value = 68865
$ws.Range("Z76").Value = 1405
$ws.Range("AA76").Value = 0
$ws.Range("AB76").Value = 1405
$ws.Range("AC76").Value = 13564
$ws.Range("AD76").Value = 9056
$ws.Range("AE76").Value = 3625
$ws.Range("AF76").Value = 5431
$ws.Range("AG76").Value = 37
$ws.Range("AH76").Value = 9
$ws.Range("AI76").Value = 28
